$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new concert row (row 16)
$ws.Range("B16").Value = "Melegnano (MI)"
$ws.Range("C16").Value = "Musicolepsia"
$ws.Range("A16").Value = "Plug & Play Jam Session (Apr 6, 2023)"
$ws.Range("D16").Value = 45.3579709544833
$ws.Range("E16").Value = 9.31467101016586
$ws.Range("F16").Value = 2023
$ws.Range("G16").Value = "06/04/2023"
$ws.Range("H16").Value = '<iframe width="300" height="169" src="https://www.youtube.com/embed/playlist?list=PLhIw1_0YGPEStVIUkVyv2ZB4PlUeK02QW"></iframe>'

# Match the date format style used in column G (numFmtId 49 -> "@")
$ws.Range("G16").NumberFormat = "@"

# These two playlist cells carried a leftover "applyFill" style that's no
# longer used anywhere else in the sheet; restore them to the default style.
$ws.Range("H12:H13").Style = "Normal"

# Select the new row's first cell, like the author did
$ws.Range("A16").Select()
